$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112228058
$ws.Range("B2").Value = 90466
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 4769
$ws.Range("F2").Value = "Svavelriska"
$ws.Range("G2").Value = "Lactarius scrobiculatus"
$ws.Range("H2").Value = "(Scop.:Fr.) Fr."
$ws.Range("P2").Value = "Nils-Andersknulen (Nils-Andersknulen), Jmt"
$ws.Range("Q2").Value = 496258
$ws.Range("R2").Value = 6934460
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = "Västernorrland"
$ws.Range("U2").Value = "Ånge"
$ws.Range("V2").Value = "Jämtland"
$ws.Range("W2").Value = "Haverö"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2023-09-21"
$ws.Range("Y2").Style = "Normal"
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "11:47"
$ws.Range("Z2").Style = "Normal"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2023-09-21"
$ws.Range("AA2").Style = "Normal"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "11:47"
$ws.Range("AB2").Style = "Normal"
$ws.Range("AD2").Value = $False
$ws.Range("AE2").Value = $False
$ws.Range("AG2").Value = $False
$ws.Range("AW2").Value = "Håkan Blomqvist"
$ws.Range("AX2").Value = "Håkan Blomqvist"

# Row 3
$ws.Range("A3").Value = 112227657
$ws.Range("B3").Value = 81371
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 1312
$ws.Range("F3").Value = "Gammelgransskål"
$ws.Range("G3").Value = "Pseudographis pinicola"
$ws.Range("H3").Value = "(Nyl.) Rehm"
$ws.Range("P3").Value = "Nils-Andersknulen (Nils-Andersknulen), Jmt"
$ws.Range("Q3").Value = 496238
$ws.Range("R3").Value = 6934504
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Västernorrland"
$ws.Range("U3").Value = "Ånge"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Haverö"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-09-21"
$ws.Range("Y3").Style = "Normal"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "11:20"
$ws.Range("Z3").Style = "Normal"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-09-21"
$ws.Range("AA3").Style = "Normal"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "11:20"
$ws.Range("AB3").Style = "Normal"
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AG3").Value = $False
$ws.Range("AW3").Value = "Håkan Blomqvist"
$ws.Range("AX3").Value = "Håkan Blomqvist"

# Row 4
$ws.Range("A4").Value = 112227891
$ws.Range("B4").Value = 89535
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1108
$ws.Range("F4").Value = "Harticka"
$ws.Range("G4").Value = "Pelloporus leporinus"
$ws.Range("H4").Value = "(Fr.) Krieglst."
$ws.Range("P4").Value = "Nils-Andersknulen (Nils-Andersknulen), Jmt"
$ws.Range("Q4").Value = 496245
$ws.Range("R4").Value = 6934459
$ws.Range("S4").Value = 25
$ws.Range("T4").Value = "Västernorrland"
$ws.Range("U4").Value = "Ånge"
$ws.Range("V4").Value = "Jämtland"
$ws.Range("W4").Value = "Haverö"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-09-21"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "11:32"
$ws.Range("Z4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-09-21"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "11:32"
$ws.Range("AB4").Style = "Normal"
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AG4").Value = $False
$ws.Range("AW4").Value = "Håkan Blomqvist"
$ws.Range("AX4").Value = "Håkan Blomqvist"

# Row 5
$ws.Range("A5").Value = 112228055
$ws.Range("B5").Value = 89503
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 5447
$ws.Range("F5").Value = "Vedticka"
$ws.Range("G5").Value = "Fuscoporia viticola"
$ws.Range("H5").Value = "(Schwein.) Murrill"
$ws.Range("P5").Value = "Nils-Andersknulen (Nils-Andersknulen), Jmt"
$ws.Range("Q5").Value = 496258
$ws.Range("R5").Value = 6934460
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = "Västernorrland"
$ws.Range("U5").Value = "Ånge"
$ws.Range("V5").Value = "Jämtland"
$ws.Range("W5").Value = "Haverö"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-21"
$ws.Range("Y5").Style = "Normal"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "11:46"
$ws.Range("Z5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-21"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "11:46"
$ws.Range("AB5").Style = "Normal"
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AG5").Value = $False
$ws.Range("AW5").Value = "Håkan Blomqvist"
$ws.Range("AX5").Value = "Håkan Blomqvist"

# Row 6
$ws.Range("A6").Value = 112228190
$ws.Range("B6").Value = 89557
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma"
$ws.Range("H6").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P6").Value = "Nils-Andersknulen (Nils-Andersknulen), Jmt"
$ws.Range("Q6").Value = 496305
$ws.Range("R6").Value = 6934462
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = "Västernorrland"
$ws.Range("U6").Value = "Ånge"
$ws.Range("V6").Value = "Jämtland"
$ws.Range("W6").Value = "Haverö"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-21"
$ws.Range("Y6").Style = "Normal"
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "11:47"
$ws.Range("Z6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-21"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "11:47"
$ws.Range("AB6").Style = "Normal"
$ws.Range("AD6").Value = $False
$ws.Range("AE6").Value = $False
$ws.Range("AG6").Value = $False
$ws.Range("AW6").Value = "Håkan Blomqvist"
$ws.Range("AX6").Value = "Håkan Blomqvist"

# Row 7
$ws.Range("A7").Value = 112228201
$ws.Range("B7").Value = 89539
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = "Ullticka"
$ws.Range("G7").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P7").Value = "Nils-Andersknulen (Nils-Andersknulen), Jmt"
$ws.Range("Q7").Value = 496302
$ws.Range("R7").Value = 6934437
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Västernorrland"
$ws.Range("U7").Value = "Ånge"
$ws.Range("V7").Value = "Jämtland"
$ws.Range("W7").Value = "Haverö"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-09-21"
$ws.Range("Y7").Style = "Normal"
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = "11:57"
$ws.Range("Z7").Style = "Normal"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-09-21"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = "11:57"
$ws.Range("AB7").Style = "Normal"
$ws.Range("AD7").Value = $False
$ws.Range("AE7").Value = $False
$ws.Range("AG7").Value = $False
$ws.Range("AW7").Value = "Håkan Blomqvist"
$ws.Range("AX7").Value = "Håkan Blomqvist"

# Row 8
$ws.Range("A8").Value = 112305970
$ws.Range("B8").Value = 90199
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 898
$ws.Range("F8").Value = "Blackticka"
$ws.Range("G8").Value = "Steccherinum collabens"
$ws.Range("H8").Value = "(Fr.) Vesterholt"
$ws.Range("P8").Value = "NilsAndersknulen, Jmt"
$ws.Range("Q8").Value = 496348
$ws.Range("R8").Value = 6934464
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = "Västernorrland"
$ws.Range("U8").Value = "Ånge"
$ws.Range("V8").Value = "Jämtland"
$ws.Range("W8").Value = "Haverö"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-09-21"
$ws.Range("Y8").Style = "Normal"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-09-21"
$ws.Range("AA8").Style = "Normal"
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AG8").Value = $False
$ws.Range("AW8").Value = "Håkan Blomqvist"
$ws.Range("AX8").Value = "Håkan Blomqvist"

# Row 8 empty placeholder cells
$ws.Range("I8").Value = "'"
$ws.Range("I8").Style = "Normal"
$ws.Range("J8").Value = "'"
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").Value = "'"
$ws.Range("K8").Style = "Normal"
$ws.Range("N8").Value = "'"
$ws.Range("N8").Style = "Normal"
$ws.Range("AF8").Value = "'"
$ws.Range("AF8").Style = "Normal"
$ws.Range("AT8").Value = "'"
$ws.Range("AT8").Style = "Normal"
$ws.Range("AY8").Value = "'"
$ws.Range("AY8").Style = "Normal"
